$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old cells that held the text (rows 13 and 14, column A)
$ws.Range("A13").Value = $null
$ws.Range("A14").Value = $null

# Add the new cells H2 and H3 with the same text values
$ws.Range("H2").Value = "Number of iterations: 7500"
$ws.Range("H3").Value = "Mutation Probability: 0.55"

# Update the selection to H2:H3 with active cell H2
$ws.Range("H2:H3").Select()

# Update window position (best-effort; window geometry is session/runtime
# metadata and may not be persisted by every host, but we set it anyway so
# that a host which does track it reflects the new position).
$excel.ActiveWindow.Left = 32160
$excel.ActiveWindow.Top = 4215
